$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Trimestre (quarter) date column C for rows 2-10.
# Force text format first so Excel keeps the literal string instead of
# auto-converting it to a date serial number, then restore the "Normal"
# style so no stray number-format styling is left on the cells.
$ws.Range("C2:C10").NumberFormat = "@"
$ws.Range("C2:C10").Value = "01/10/2023"
$ws.Range("C2:C10").Style = "Normal"

# Update the Valor (value) column D with new figures
$ws.Range("D2").Value = 56.27161270040868
$ws.Range("D3").Value = 55.74315341664451
$ws.Range("D4").Value = 55.47008366878122
$ws.Range("D5").Value = 54.92521223554777
$ws.Range("D6").Value = 54.34933287125282
$ws.Range("D7").Value = 53.89736477115118
$ws.Range("D8").Value = 44.94334872010072
$ws.Range("D9").Value = 43.34118375135377
$ws.Range("D10").Value = 50.46478162897979

# Row 7 region name changed from "Mato Grosso do Sul" to "Mato Grosso"
$ws.Range("A7").Value = "Mato Grosso"

# Row 8 ranking changed from 18º to 19º
$ws.Range("E8").Value = "19º"
